# Add data for 2022-06-01: update "through" date from 05-23 to 05-24,
# rename sheet, and bump the May / Total row values for 2017-2022.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab / sheet name
$ws.Name = "Through 2022-05-24"

# Update the month label for row 6 (May)
$ws.Range("A6").Value = "May (through 05-24)"

# Row 6 (May) updated figures for columns D..I (2017-2022)
$ws.Range("D6").Value = 47
$ws.Range("E6").Value = 39
$ws.Range("G6").Value = 46
$ws.Range("H6").Value = 93
$ws.Range("I6").Value = 89

# Row 7 (Total) updated figures for columns D..I (2017-2022)
$ws.Range("D7").Value = 300
$ws.Range("E7").Value = 285
$ws.Range("G7").Value = 308
$ws.Range("H7").Value = 616
$ws.Range("I7").Value = 641
